$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the stray Google hyperlink (and its displayed text) that lived in
# F103, while leaving the cell's Hyperlink style intact.
$ws.Range("F103").Hyperlinks.Delete()
$ws.Range("F103").Value = ""

# Remove the three rows whose rules are no longer part of the list:
#   row 8  -> squid:S2076 "Values passed to OS commands should be sanitized"
#   row 10 -> squid:S2078 "Values passed to LDAP queries should be sanitized"
#   row 26 -> squid:S3318 "Untrusted data should not be stored in sessions"
# Delete from the bottom up so earlier row numbers stay valid.
$ws.Rows(26).Delete()
$ws.Rows(10).Delete()
$ws.Rows(8).Delete()

# Match the saved view/selection state from the edited workbook.
$ws.Range("F100").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 73
$win.ScrollColumn = 1
